# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new row 312 and 313) for "Terminal
# Hortofrutícola Agro Chillán - Mandarina" and push the existing data
# down by two rows (old row 312 -> new row 314, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 312; everything from 312 on shifts down to 314+
$ws.Rows.Item(312).Resize(2).Insert()

# ---- New row 312 ----
$ws.Range("A312").Value = 7
$ws.Range("B312").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C312").Value = "Ñuble"
$ws.Range("D312").Value = 45132
$ws.Range("E312").Value = 16
$ws.Range("F312").Value = "Fruta"
$ws.Range("G312").Value = 100102
$ws.Range("H312").Value = "Cítricos"
$ws.Range("I312").Value = 100102004
$ws.Range("J312").Value = "Mandarina"
$ws.Range("K312").Value = "Clementina"
$ws.Range("L312").Value = "Primera"
$ws.Range("M312").Value = 150
$ws.Range("N312").Value = 9000
$ws.Range("O312").Value = 9000
$ws.Range("P312").Value = 9000
$ws.Range("Q312").Value = "$/bandeja 10 kilos"
$ws.Range("R312").Value = "Región de O'Higgins"
$ws.Range("S312").Value = 900
$ws.Range("T312").Value = 10

# ---- New row 313 ----
$ws.Range("A313").Value = 7
$ws.Range("B313").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C313").Value = "Ñuble"
$ws.Range("D313").Value = 45132
$ws.Range("E313").Value = 16
$ws.Range("F313").Value = "Fruta"
$ws.Range("G313").Value = 100102
$ws.Range("H313").Value = "Cítricos"
$ws.Range("I313").Value = 100102004
$ws.Range("J313").Value = "Mandarina"
$ws.Range("K313").Value = "Clementina"
$ws.Range("L313").Value = "Segunda"
$ws.Range("M313").Value = 150
$ws.Range("N313").Value = 8000
$ws.Range("O313").Value = 8000
$ws.Range("P313").Value = 8000
$ws.Range("Q313").Value = "$/bandeja 10 kilos"
$ws.Range("R313").Value = "Región de O'Higgins"
$ws.Range("S313").Value = 800
$ws.Range("T313").Value = 10
